$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$wsAssets = $wb.Worksheets.Item("Assets")

# --- Remove hyperlink formatting from B13 / B16 / (old) B24 ---
# Deleting hyperlinks on the sheet also drops the <hyperlinks> element
# and its relationships entirely (matches target removing all 3 links).
$ws.Range("B13").Hyperlinks.Delete()
$ws.Range("B13").Style = "Normal"
$ws.Range("B16").Style = "Normal"
$ws.Range("B24").Style = "Normal"

# Drop the now-unused "Hyperlink" cell style/font definitions.
$wb.Styles.Item("Hyperlink").Delete()

# --- Rearrange / extend the Settings rows 21-32 ---
$ws.Range("A21:C32").ClearContents()
$ws.Range("A21:C32").Style = "Normal"

$ws.Range("A21").Value = "AttachmentName"
$ws.Range("B21").Value = "Colegi.xlsx"
$ws.Range("C21").Value = "The name of attachment to be looked after."

$ws.Range("A22").Value = "AttachmentSheetName"
$ws.Range("B22").Value = "Sheet1"
$ws.Range("C22").Value = "Name of used sheet in excel file"

$ws.Range("A24").Value = "OutputReportPath"
$ws.Range("B24").Value = "C:\Users\DariusDangi\Desktop\OutputReport.xlsx"
$ws.Range("C24").Value = "Path to where the output reports its saved. NOTE: change this acording to your case."

$ws.Range("A25").Value = "OutputReportSheetName"
$ws.Range("B25").Value = "Sheet1"
$ws.Range("C25").Value = "The name of excel sheet that is used when it's created."

$ws.Range("A26").Value = "OutputReportColumns"
$ws.Range("B26").Value = "Name,Email,Animal,Status"
$ws.Range("C26").Value = "Column names that are used to create the output report."

$ws.Range("A28").Value = "OutlookInputEmailAddress"
$ws.Range("B28").Value = "darius.dangi@fwfcompany.com"
$ws.Range("C28").Value = "Email Address used to send mail messages"

$ws.Range("A30").Value = "DictionaryTransactionItemsColumns"
$ws.Range("B30").Value = "Name,Email,Animal "
$ws.Range("C30").Value = "These are the keys that are going to be used in the queue items."

$ws.Range("A31").Value = "SearchingForColumn"
$ws.Range("B31").Value = "Animal "
$ws.Range("C31").Value = "We are going to search this column to see if it contains dog/cat."

$ws.Range("A32").Value = "DictionarySenderEmailColumn"
$ws.Range("B32").Value = "SenderEmail"
$ws.Range("C32").Value = "This is not part of transactionItem, so we have to use another value."

# --- View state tweaks ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C16").Select()

$wsAssets.Application.ActiveWindow.ScrollRow = 1
